# "brower config and commonDataProvider"
# - Rename the existing sheet to "invalidCredentialTest"
# - Add a new "validCredentialTest" sheet (common data provider for valid creds)
#   after it, with its own username/password/language/expectedTitle rows
# - Move the selection on the original sheet, and select+activate the new sheet

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "invalidCredentialTest"

# Insert the new sheet right after the (renamed) first sheet.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "validCredentialTest"

# Header row, shared with the invalidCredentialTest sheet.
$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "password"
$ws2.Range("C1").Value = "language"
$ws2.Range("D1").Value = "expectedTitle"

# Fill column-by-column so new shared-string entries are interned in the
# same order the authored workbook used (admin, physician, pass, OpenEMR).
$ws2.Range("A2").Value = "admin"
$ws2.Range("A3").Value = "physician"

$ws2.Range("B2").Value = "pass"
$ws2.Range("B3").Value = "physician"

$ws2.Range("C2").Value = "English (Indian)"
$ws2.Range("C3").Value = "Dutch"

$ws2.Range("D2").Value = "OpenEMR"
$ws2.Range("D3").Value = "OpenEMR"

# Leave the old sheet's selection where the author left it.
$ws1.Range("G16").Select() | Out-Null

# The new sheet ends up the active / selected tab.
$ws2.Activate()
$ws2.Range("D5").Select() | Out-Null
